$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing age-range labels
$ws.Range("A2").Value = "10-15 years"
$ws.Range("A4").Value = "20-25 years"

# Zero-out the existing Female/Male counts for rows 3 and 4
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0

# New age-range rows (5-14), each with style from the A-column label cells
# and zeroed Female/Male counts, copying the format used by the existing
# label cells (A2:A4).
$ageRanges = @(
  "25-30 years",
  "30-35 years",
  "35-40 years",
  "40-45 years",
  "45-50 years",
  "50-55 years",
  "55-60 years",
  "60-65 years",
  "70-75 years",
  "75-80 years"
)

$row = 5
foreach ($age in $ageRanges) {
  $ws.Cells.Item(2, 1).Copy()
  $ws.Cells.Item($row, 1).PasteSpecial(-4122)
  $ws.Cells.Item($row, 1).Value = $age
  $ws.Cells.Item($row, 2).Value = 0
  $ws.Cells.Item($row, 3).Value = 0
  $row = $row + 1
}
